$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "apelido" (nickname) value on row 2 from Wilkerbn007 -> Wilkerbn010
$ws.Range("A2").Value = "Wilkerbn010"

# Add four new header columns (M1:P1)
$ws.Range("M1").Value = "campoCategoria(Pesquisa)"
$ws.Range("N1").Value = "campoProduto (Pesquisa)"
$ws.Range("O1").Value = "campoCategoria(Tela Principal)"
$ws.Range("P1").Value = "campoProduto (Tela Principal)"

# Add the corresponding data row (M2:P2)
$ws.Range("M2").Value = "Laptop"
$ws.Range("N2").Value = "HP PAVILION 15T TOUCH LAPTOP"
$ws.Range("O2").Value = "LAPTOPS"
$ws.Range("P2").Value = "HP Pavilion 15z Laptop"

# Match the saved view state: scrolled to show column M, selection on P4
$ws.Application.ActiveWindow.ScrollColumn = 13
$ws.Range("P4").Select()
